$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: the "Find One" subtask is dropped, the ideation date moves out
#     to 6/22/2019, and a new comment/hyperlink about the matplotlib guide's
#     reader comments section is recorded.
$ws.Range("D16").ClearContents()
$ws.Range("C16").Value = "6/22/2019"
$ws.Hyperlinks.Add($ws.Range("G16"), "https://realpython.com/python-matplotlib-guide/", "reader-comments", "", "https://realpython.com/python-matplotlib-guide/#reader-comments") | Out-Null
$ws.Range("G16").Style = "Hyperlink"

# --- Insert a new blank separator row under it (everything below shifts down by one).
$ws.Rows("17").Insert()

# --- Tidy the "Create Quad model..." task text and the trailing-space "Ongoing " status.
# (after the row insert above, these now live at row 20 instead of the old row 19)
$ws.Range("B20").Value = "Create Quad model for forces and torques"
$ws.Range("F20").Value = "Ongoing"

# --- Update the view: scrolled/selected cell moves from G19 to F20.
$excel.Goto($ws.Range("A13"), $true) | Out-Null
$ws.Range("F20").Select() | Out-Null
